$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 and 46: swap EnergySwap/Mantle entries with updated values
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.597"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.37%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.50"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.04%  "

# Remaining price/volume updates across the coin list
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.638.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.626.85'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.03'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.576'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.633.50'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.083.24'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.689.60'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.73'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.24%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.627.49'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '345.87'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.44'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.18'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.13'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.57'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.08%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0801'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.84'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.32'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.972'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.97'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.14'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.61'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.15%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.63'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '280.06'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0983'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.36%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.31'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.979.80'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.68%  '
